# Update Num_Inclusions (column C) values for several rows
# per the re-analysis described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C32").Value = 12
$ws.Range("C34").Value = 2
$ws.Range("C42").Value = 8
$ws.Range("C72").Value = 4
$ws.Range("C73").Value = 25
$ws.Range("C84").Value = 1
$ws.Range("C93").Value = 2
$ws.Range("C95").Value = 5
$ws.Range("C175").Value = 1
